$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = 2
$ws1.Range("B2").Value = 1
$ws1.Activate()
$ws1.Range("M5").Select()
